# Apply the edits described by the diff:
# - Update several values in column B (Goal1 values)
# - Move the active/selected cell to D12

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values
$ws.Range("B2").Value = 37.39
$ws.Range("B3").Value = 69.49
$ws.Range("B5").Value = 0.627
$ws.Range("B6").Value = 0.327
$ws.Range("B7").Value = 0.327
$ws.Range("B8").Value = 0.627

# Move the active selection to D12
$ws.Range("D12").Select()
